# Auto-generated: rewrite rows 525-633 with updated weekly Limon price data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44641
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 160
$arr[0,13] = 21000
$arr[0,14] = 22000
$arr[0,15] = 21500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1344
$arr[0,19] = 16
$ws.Range("A525:T525").Value = $arr
$ws.Range("D525").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44641
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 120
$arr[0,13] = 19000
$arr[0,14] = 20000
$arr[0,15] = 19500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1219
$arr[0,19] = 16
$ws.Range("A526:T526").Value = $arr
$ws.Range("D526").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44421
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 300
$arr[0,13] = 3800
$arr[0,14] = 4000
$arr[0,15] = 3900
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 244
$arr[0,19] = 16
$ws.Range("A527:T527").Value = $arr
$ws.Range("D527").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44421
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 220
$arr[0,13] = 3000
$arr[0,14] = 3500
$arr[0,15] = 3227
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 202
$arr[0,19] = 16
$ws.Range("A528:T528").Value = $arr
$ws.Range("D528").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44421
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "3a amarillo"
$arr[0,12] = 100
$arr[0,13] = 2600
$arr[0,14] = 2800
$arr[0,15] = 2700
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 169
$arr[0,19] = 16
$ws.Range("A529:T529").Value = $arr
$ws.Range("D529").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44329
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 120
$arr[0,13] = 14000
$arr[0,14] = 15000
$arr[0,15] = 14500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 906
$arr[0,19] = 16
$ws.Range("A530:T530").Value = $arr
$ws.Range("D530").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44329
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 12000
$arr[0,14] = 13000
$arr[0,15] = 12500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 781
$arr[0,19] = 16
$ws.Range("A531:T531").Value = $arr
$ws.Range("D531").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44637
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 60
$arr[0,13] = 20000
$arr[0,14] = 20000
$arr[0,15] = 20000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1250
$arr[0,19] = 16
$ws.Range("A532:T532").Value = $arr
$ws.Range("D532").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44637
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 120
$arr[0,13] = 20000
$arr[0,14] = 21000
$arr[0,15] = 20500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1281
$arr[0,19] = 16
$ws.Range("A533:T533").Value = $arr
$ws.Range("D533").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44208
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 170
$arr[0,13] = 19000
$arr[0,14] = 20000
$arr[0,15] = 19471
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región Metropolitana"
$arr[0,18] = 1217
$arr[0,19] = 16
$ws.Range("A534:T534").Value = $arr
$ws.Range("D534").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44208
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 220
$arr[0,13] = 21000
$arr[0,14] = 22000
$arr[0,15] = 21636
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de Coquimbo"
$arr[0,18] = 1352
$arr[0,19] = 16
$ws.Range("A535:T535").Value = $arr
$ws.Range("D535").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44445
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 200
$arr[0,13] = 3800
$arr[0,14] = 4000
$arr[0,15] = 3900
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 244
$arr[0,19] = 16
$ws.Range("A536:T536").Value = $arr
$ws.Range("D536").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44445
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 200
$arr[0,13] = 3000
$arr[0,14] = 3500
$arr[0,15] = 3250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 203
$arr[0,19] = 16
$ws.Range("A537:T537").Value = $arr
$ws.Range("D537").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44445
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "3a amarillo"
$arr[0,12] = 160
$arr[0,13] = 2600
$arr[0,14] = 2800
$arr[0,15] = 2700
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 169
$arr[0,19] = 16
$ws.Range("A538:T538").Value = $arr
$ws.Range("D538").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44524
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 7000
$arr[0,14] = 7500
$arr[0,15] = 7250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 453
$arr[0,19] = 16
$ws.Range("A539:T539").Value = $arr
$ws.Range("D539").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44355
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 120
$arr[0,13] = 9000
$arr[0,14] = 10000
$arr[0,15] = 9500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 594
$arr[0,19] = 16
$ws.Range("A540:T540").Value = $arr
$ws.Range("D540").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44355
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 7000
$arr[0,14] = 8000
$arr[0,15] = 7500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 469
$arr[0,19] = 16
$ws.Range("A541:T541").Value = $arr
$ws.Range("D541").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44530
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 7500
$arr[0,14] = 8000
$arr[0,15] = 7750
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 484
$arr[0,19] = 16
$ws.Range("A542:T542").Value = $arr
$ws.Range("D542").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44530
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 160
$arr[0,13] = 6500
$arr[0,14] = 7000
$arr[0,15] = 6750
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 422
$arr[0,19] = 16
$ws.Range("A543:T543").Value = $arr
$ws.Range("D543").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44483
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 240
$arr[0,13] = 4000
$arr[0,14] = 4500
$arr[0,15] = 4250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 266
$arr[0,19] = 16
$ws.Range("A544:T544").Value = $arr
$ws.Range("D544").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44483
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 240
$arr[0,13] = 3000
$arr[0,14] = 3500
$arr[0,15] = 3250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 203
$arr[0,19] = 16
$ws.Range("A545:T545").Value = $arr
$ws.Range("D545").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44483
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "3a amarillo"
$arr[0,12] = 60
$arr[0,13] = 2800
$arr[0,14] = 2800
$arr[0,15] = 2800
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 175
$arr[0,19] = 16
$ws.Range("A546:T546").Value = $arr
$ws.Range("D546").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44294
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 15000
$arr[0,14] = 16000
$arr[0,15] = 15500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 969
$arr[0,19] = 16
$ws.Range("A547:T547").Value = $arr
$ws.Range("D547").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44294
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 13000
$arr[0,14] = 14000
$arr[0,15] = 13500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 844
$arr[0,19] = 16
$ws.Range("A548:T548").Value = $arr
$ws.Range("D548").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44617
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 120
$arr[0,13] = 18000
$arr[0,14] = 19000
$arr[0,15] = 18500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1156
$arr[0,19] = 16
$ws.Range("A549:T549").Value = $arr
$ws.Range("D549").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44557
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 80
$arr[0,13] = 13500
$arr[0,14] = 13500
$arr[0,15] = 13500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 844
$arr[0,19] = 16
$ws.Range("A550:T550").Value = $arr
$ws.Range("D550").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44489
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 5000
$arr[0,14] = 5500
$arr[0,15] = 5250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 328
$arr[0,19] = 16
$ws.Range("A551:T551").Value = $arr
$ws.Range("D551").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44489
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 160
$arr[0,13] = 4000
$arr[0,14] = 4500
$arr[0,15] = 4250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 266
$arr[0,19] = 16
$ws.Range("A552:T552").Value = $arr
$ws.Range("D552").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44264
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 160
$arr[0,13] = 17000
$arr[0,14] = 18000
$arr[0,15] = 17500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1094
$arr[0,19] = 16
$ws.Range("A553:T553").Value = $arr
$ws.Range("D553").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44264
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 80
$arr[0,13] = 16000
$arr[0,14] = 16000
$arr[0,15] = 16000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1000
$arr[0,19] = 16
$ws.Range("A554:T554").Value = $arr
$ws.Range("D554").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44396
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 4000
$arr[0,14] = 4500
$arr[0,15] = 4250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 266
$arr[0,19] = 16
$ws.Range("A555:T555").Value = $arr
$ws.Range("D555").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44396
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 3000
$arr[0,14] = 3500
$arr[0,15] = 3250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 203
$arr[0,19] = 16
$ws.Range("A556:T556").Value = $arr
$ws.Range("D556").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44232
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 220
$arr[0,13] = 18000
$arr[0,14] = 18500
$arr[0,15] = 18227
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Quillota"
$arr[0,18] = 1139
$arr[0,19] = 16
$ws.Range("A557:T557").Value = $arr
$ws.Range("D557").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44232
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 280
$arr[0,13] = 18000
$arr[0,14] = 19000
$arr[0,15] = 18536
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de Coquimbo"
$arr[0,18] = 1158
$arr[0,19] = 16
$ws.Range("A558:T558").Value = $arr
$ws.Range("D558").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44279
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 17000
$arr[0,14] = 18000
$arr[0,15] = 17500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1094
$arr[0,19] = 16
$ws.Range("A559:T559").Value = $arr
$ws.Range("D559").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44279
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 120
$arr[0,13] = 17000
$arr[0,14] = 18000
$arr[0,15] = 17500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1094
$arr[0,19] = 16
$ws.Range("A560:T560").Value = $arr
$ws.Range("D560").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44504
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 6500
$arr[0,14] = 6800
$arr[0,15] = 6650
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 416
$arr[0,19] = 16
$ws.Range("A561:T561").Value = $arr
$ws.Range("D561").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44504
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 5500
$arr[0,14] = 6000
$arr[0,15] = 5750
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 359
$arr[0,19] = 16
$ws.Range("A562:T562").Value = $arr
$ws.Range("D562").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44572
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 200
$arr[0,13] = 14000
$arr[0,14] = 15000
$arr[0,15] = 14500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 906
$arr[0,19] = 16
$ws.Range("A563:T563").Value = $arr
$ws.Range("D563").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44572
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 160
$arr[0,13] = 12000
$arr[0,14] = 13000
$arr[0,15] = 12500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 781
$arr[0,19] = 16
$ws.Range("A564:T564").Value = $arr
$ws.Range("D564").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44257
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 160
$arr[0,13] = 17000
$arr[0,14] = 18000
$arr[0,15] = 17500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1094
$arr[0,19] = 16
$ws.Range("A565:T565").Value = $arr
$ws.Range("D565").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44257
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 60
$arr[0,13] = 16000
$arr[0,14] = 16000
$arr[0,15] = 16000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1000
$arr[0,19] = 16
$ws.Range("A566:T566").Value = $arr
$ws.Range("D566").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44301
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 120
$arr[0,13] = 15000
$arr[0,14] = 16000
$arr[0,15] = 15500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 969
$arr[0,19] = 16
$ws.Range("A567:T567").Value = $arr
$ws.Range("D567").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44301
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 120
$arr[0,13] = 13000
$arr[0,14] = 14000
$arr[0,15] = 13500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 844
$arr[0,19] = 16
$ws.Range("A568:T568").Value = $arr
$ws.Range("D568").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44370
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 6000
$arr[0,14] = 6500
$arr[0,15] = 6250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 391
$arr[0,19] = 16
$ws.Range("A569:T569").Value = $arr
$ws.Range("D569").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44370
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 5000
$arr[0,14] = 5500
$arr[0,15] = 5250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 328
$arr[0,19] = 16
$ws.Range("A570:T570").Value = $arr
$ws.Range("D570").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44487
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 4000
$arr[0,14] = 4500
$arr[0,15] = 4250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 266
$arr[0,19] = 16
$ws.Range("A571:T571").Value = $arr
$ws.Range("D571").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44487
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 160
$arr[0,13] = 3500
$arr[0,14] = 3800
$arr[0,15] = 3650
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 228
$arr[0,19] = 16
$ws.Range("A572:T572").Value = $arr
$ws.Range("D572").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44174
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 250
$arr[0,13] = 11000
$arr[0,14] = 12000
$arr[0,15] = 11520
$arr[0,16] = "`$/malla 18 kilos"
$arr[0,17] = "Región de Coquimbo"
$arr[0,18] = 640
$arr[0,19] = 18
$ws.Range("A573:T573").Value = $arr
$ws.Range("D573").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44200
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 120
$arr[0,13] = 18000
$arr[0,14] = 19000
$arr[0,15] = 18500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 1156
$arr[0,19] = 16
$ws.Range("A574:T574").Value = $arr
$ws.Range("D574").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44200
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 60
$arr[0,13] = 21000
$arr[0,14] = 22000
$arr[0,15] = 21500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Quillota"
$arr[0,18] = 1344
$arr[0,19] = 16
$ws.Range("A575:T575").Value = $arr
$ws.Range("D575").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44200
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 16500
$arr[0,14] = 17000
$arr[0,15] = 16750
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 1047
$arr[0,19] = 16
$ws.Range("A576:T576").Value = $arr
$ws.Range("D576").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44385
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 4500
$arr[0,14] = 5000
$arr[0,15] = 4750
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 297
$arr[0,19] = 16
$ws.Range("A577:T577").Value = $arr
$ws.Range("D577").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44385
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 3400
$arr[0,14] = 3800
$arr[0,15] = 3600
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 225
$arr[0,19] = 16
$ws.Range("A578:T578").Value = $arr
$ws.Range("D578").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44236
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 120
$arr[0,13] = 18000
$arr[0,14] = 19000
$arr[0,15] = 18500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1156
$arr[0,19] = 16
$ws.Range("A579:T579").Value = $arr
$ws.Range("D579").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44236
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 120
$arr[0,13] = 16000
$arr[0,14] = 17000
$arr[0,15] = 16500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1031
$arr[0,19] = 16
$ws.Range("A580:T580").Value = $arr
$ws.Range("D580").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44221
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 240
$arr[0,13] = 19000
$arr[0,14] = 20000
$arr[0,15] = 19500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 1219
$arr[0,19] = 16
$ws.Range("A581:T581").Value = $arr
$ws.Range("D581").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44221
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 240
$arr[0,13] = 17000
$arr[0,14] = 18000
$arr[0,15] = 17500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 1094
$arr[0,19] = 16
$ws.Range("A582:T582").Value = $arr
$ws.Range("D582").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44413
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 240
$arr[0,13] = 3800
$arr[0,14] = 4000
$arr[0,15] = 3900
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 244
$arr[0,19] = 16
$ws.Range("A583:T583").Value = $arr
$ws.Range("D583").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44413
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 200
$arr[0,13] = 3000
$arr[0,14] = 3200
$arr[0,15] = 3100
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 194
$arr[0,19] = 16
$ws.Range("A584:T584").Value = $arr
$ws.Range("D584").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44413
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "3a amarillo"
$arr[0,12] = 160
$arr[0,13] = 2600
$arr[0,14] = 2800
$arr[0,15] = 2700
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 169
$arr[0,19] = 16
$ws.Range("A585:T585").Value = $arr
$ws.Range("D585").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44272
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 160
$arr[0,13] = 17000
$arr[0,14] = 18000
$arr[0,15] = 17500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1094
$arr[0,19] = 16
$ws.Range("A586:T586").Value = $arr
$ws.Range("D586").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44272
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 80
$arr[0,13] = 16000
$arr[0,14] = 16000
$arr[0,15] = 16000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1000
$arr[0,19] = 16
$ws.Range("A587:T587").Value = $arr
$ws.Range("D587").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44229
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 280
$arr[0,13] = 16500
$arr[0,14] = 17000
$arr[0,15] = 16732
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de Coquimbo"
$arr[0,18] = 1046
$arr[0,19] = 16
$ws.Range("A588:T588").Value = $arr
$ws.Range("D588").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44214
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 150
$arr[0,13] = 19500
$arr[0,14] = 20000
$arr[0,15] = 19783
$arr[0,16] = "`$/bandeja 15 kilos"
$arr[0,17] = "Provincia de Quillota"
$arr[0,18] = 1319
$arr[0,19] = 15
$ws.Range("A589:T589").Value = $arr
$ws.Range("D589").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44214
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 160
$arr[0,13] = 21000
$arr[0,14] = 22000
$arr[0,15] = 21531
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de Coquimbo"
$arr[0,18] = 1346
$arr[0,19] = 16
$ws.Range("A590:T590").Value = $arr
$ws.Range("D590").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44299
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 160
$arr[0,13] = 16000
$arr[0,14] = 17000
$arr[0,15] = 16500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1031
$arr[0,19] = 16
$ws.Range("A591:T591").Value = $arr
$ws.Range("D591").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44299
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 80
$arr[0,13] = 15000
$arr[0,14] = 15000
$arr[0,15] = 15000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 938
$arr[0,19] = 16
$ws.Range("A592:T592").Value = $arr
$ws.Range("D592").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44610
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 120
$arr[0,13] = 16000
$arr[0,14] = 17000
$arr[0,15] = 16500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1031
$arr[0,19] = 16
$ws.Range("A593:T593").Value = $arr
$ws.Range("D593").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44610
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 100
$arr[0,13] = 14000
$arr[0,14] = 15000
$arr[0,15] = 14500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 906
$arr[0,19] = 16
$ws.Range("A594:T594").Value = $arr
$ws.Range("D594").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44312
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 15000
$arr[0,14] = 16000
$arr[0,15] = 15500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 969
$arr[0,19] = 16
$ws.Range("A595:T595").Value = $arr
$ws.Range("D595").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44312
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 80
$arr[0,13] = 14000
$arr[0,14] = 14000
$arr[0,15] = 14000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 875
$arr[0,19] = 16
$ws.Range("A596:T596").Value = $arr
$ws.Range("D596").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44399
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 120
$arr[0,13] = 4000
$arr[0,14] = 4500
$arr[0,15] = 4250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 266
$arr[0,19] = 16
$ws.Range("A597:T597").Value = $arr
$ws.Range("D597").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44399
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 3500
$arr[0,14] = 3600
$arr[0,15] = 3550
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 222
$arr[0,19] = 16
$ws.Range("A598:T598").Value = $arr
$ws.Range("D598").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44615
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 160
$arr[0,13] = 16000
$arr[0,14] = 17000
$arr[0,15] = 16500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1031
$arr[0,19] = 16
$ws.Range("A599:T599").Value = $arr
$ws.Range("D599").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44615
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 80
$arr[0,13] = 15000
$arr[0,14] = 15000
$arr[0,15] = 15000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 938
$arr[0,19] = 16
$ws.Range("A600:T600").Value = $arr
$ws.Range("D600").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44522
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 200
$arr[0,13] = 6000
$arr[0,14] = 6500
$arr[0,15] = 6250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 391
$arr[0,19] = 16
$ws.Range("A601:T601").Value = $arr
$ws.Range("D601").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44522
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 5000
$arr[0,14] = 5500
$arr[0,15] = 5250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 328
$arr[0,19] = 16
$ws.Range("A602:T602").Value = $arr
$ws.Range("D602").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44543
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 8000
$arr[0,14] = 8500
$arr[0,15] = 8250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Limarí"
$arr[0,18] = 516
$arr[0,19] = 16
$ws.Range("A603:T603").Value = $arr
$ws.Range("D603").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44543
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 160
$arr[0,13] = 10000
$arr[0,14] = 11000
$arr[0,15] = 10500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 656
$arr[0,19] = 16
$ws.Range("A604:T604").Value = $arr
$ws.Range("D604").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44543
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 160
$arr[0,13] = 7000
$arr[0,14] = 7500
$arr[0,15] = 7250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Limarí"
$arr[0,18] = 453
$arr[0,19] = 16
$ws.Range("A605:T605").Value = $arr
$ws.Range("D605").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44167
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 120
$arr[0,13] = 7000
$arr[0,14] = 7500
$arr[0,15] = 7250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 453
$arr[0,19] = 16
$ws.Range("A606:T606").Value = $arr
$ws.Range("D606").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44167
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 6000
$arr[0,14] = 6500
$arr[0,15] = 6250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 391
$arr[0,19] = 16
$ws.Range("A607:T607").Value = $arr
$ws.Range("D607").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44277
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 120
$arr[0,13] = 17000
$arr[0,14] = 18000
$arr[0,15] = 17500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1094
$arr[0,19] = 16
$ws.Range("A608:T608").Value = $arr
$ws.Range("D608").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44277
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 80
$arr[0,13] = 16000
$arr[0,14] = 16000
$arr[0,15] = 16000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1000
$arr[0,19] = 16
$ws.Range("A609:T609").Value = $arr
$ws.Range("D609").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44258
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 160
$arr[0,13] = 17000
$arr[0,14] = 18000
$arr[0,15] = 17500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1094
$arr[0,19] = 16
$ws.Range("A610:T610").Value = $arr
$ws.Range("D610").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44258
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 80
$arr[0,13] = 16000
$arr[0,14] = 16000
$arr[0,15] = 16000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1000
$arr[0,19] = 16
$ws.Range("A611:T611").Value = $arr
$ws.Range("D611").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44390
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 160
$arr[0,13] = 4500
$arr[0,14] = 5000
$arr[0,15] = 4750
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 297
$arr[0,19] = 16
$ws.Range("A612:T612").Value = $arr
$ws.Range("D612").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44349
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 120
$arr[0,13] = 9000
$arr[0,14] = 10000
$arr[0,15] = 9500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 594
$arr[0,19] = 16
$ws.Range("A613:T613").Value = $arr
$ws.Range("D613").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44349
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 7000
$arr[0,14] = 8000
$arr[0,15] = 7500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 469
$arr[0,19] = 16
$ws.Range("A614:T614").Value = $arr
$ws.Range("D614").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44498
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 300
$arr[0,13] = 5000
$arr[0,14] = 5500
$arr[0,15] = 5250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 328
$arr[0,19] = 16
$ws.Range("A615:T615").Value = $arr
$ws.Range("D615").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44498
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 160
$arr[0,13] = 4000
$arr[0,14] = 4500
$arr[0,15] = 4250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 266
$arr[0,19] = 16
$ws.Range("A616:T616").Value = $arr
$ws.Range("D616").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44179
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 240
$arr[0,13] = 12000
$arr[0,14] = 13000
$arr[0,15] = 12667
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de Coquimbo"
$arr[0,18] = 792
$arr[0,19] = 16
$ws.Range("A617:T617").Value = $arr
$ws.Range("D617").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44418
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 220
$arr[0,13] = 3800
$arr[0,14] = 4000
$arr[0,15] = 3909
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 244
$arr[0,19] = 16
$ws.Range("A618:T618").Value = $arr
$ws.Range("D618").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44418
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 160
$arr[0,13] = 3000
$arr[0,14] = 3500
$arr[0,15] = 3250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 203
$arr[0,19] = 16
$ws.Range("A619:T619").Value = $arr
$ws.Range("D619").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44418
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "3a amarillo"
$arr[0,12] = 120
$arr[0,13] = 2600
$arr[0,14] = 2800
$arr[0,15] = 2700
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 169
$arr[0,19] = 16
$ws.Range("A620:T620").Value = $arr
$ws.Range("D620").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44595
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 100
$arr[0,13] = 15000
$arr[0,14] = 15000
$arr[0,15] = 15000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 938
$arr[0,19] = 16
$ws.Range("A621:T621").Value = $arr
$ws.Range("D621").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44595
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 100
$arr[0,13] = 16000
$arr[0,14] = 16000
$arr[0,15] = 16000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1000
$arr[0,19] = 16
$ws.Range("A622:T622").Value = $arr
$ws.Range("D622").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44628
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 120
$arr[0,13] = 19000
$arr[0,14] = 20000
$arr[0,15] = 19500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1219
$arr[0,19] = 16
$ws.Range("A623:T623").Value = $arr
$ws.Range("D623").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44628
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 160
$arr[0,13] = 23000
$arr[0,14] = 24000
$arr[0,15] = 23500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Provincia de Melipilla"
$arr[0,18] = 1469
$arr[0,19] = 16
$ws.Range("A624:T624").Value = $arr
$ws.Range("D624").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44628
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 80
$arr[0,13] = 17000
$arr[0,14] = 17000
$arr[0,15] = 17000
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 1062
$arr[0,19] = 16
$ws.Range("A625:T625").Value = $arr
$ws.Range("D625").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44335
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 120
$arr[0,13] = 12000
$arr[0,14] = 13000
$arr[0,15] = 12500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 781
$arr[0,19] = 16
$ws.Range("A626:T626").Value = $arr
$ws.Range("D626").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44552
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 120
$arr[0,13] = 15000
$arr[0,14] = 16000
$arr[0,15] = 15500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 969
$arr[0,19] = 16
$ws.Range("A627:T627").Value = $arr
$ws.Range("D627").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44552
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 120
$arr[0,13] = 13000
$arr[0,14] = 14000
$arr[0,15] = 13500
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 844
$arr[0,19] = 16
$ws.Range("A628:T628").Value = $arr
$ws.Range("D628").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44544
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 120
$arr[0,13] = 12000
$arr[0,14] = 13000
$arr[0,15] = 12500
$arr[0,16] = "`$/bandeja 15 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 833
$arr[0,19] = 15
$ws.Range("A629:T629").Value = $arr
$ws.Range("D629").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44544
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a amarillo"
$arr[0,12] = 60
$arr[0,13] = 10000
$arr[0,14] = 10000
$arr[0,15] = 10000
$arr[0,16] = "`$/bandeja 15 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 667
$arr[0,19] = 15
$ws.Range("A630:T630").Value = $arr
$ws.Range("D630").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44160
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a amarillo"
$arr[0,12] = 120
$arr[0,13] = 7500
$arr[0,14] = 8000
$arr[0,15] = 7750
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 484
$arr[0,19] = 16
$ws.Range("A631:T631").Value = $arr
$ws.Range("D631").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44160
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "1a plateado"
$arr[0,12] = 160
$arr[0,13] = 8000
$arr[0,14] = 8500
$arr[0,15] = 8250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 516
$arr[0,19] = 16
$ws.Range("A632:T632").Value = $arr
$ws.Range("D632").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$arr = New-Object "object[,]" 1,20
$arr[0,0] = 7
$arr[0,1] = "Terminal Hortofrutícola Agro Chillán"
$arr[0,2] = "Ñuble"
$arr[0,3] = 44160
$arr[0,4] = 16
$arr[0,5] = "Fruta"
$arr[0,6] = 100102
$arr[0,7] = "Cítricos"
$arr[0,8] = 100102003
$arr[0,9] = "Limón"
$arr[0,10] = "Sin especificar"
$arr[0,11] = "2a plateado"
$arr[0,12] = 120
$arr[0,13] = 7000
$arr[0,14] = 7500
$arr[0,15] = 7250
$arr[0,16] = "`$/malla 16 kilos"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 453
$arr[0,19] = 16
$ws.Range("A633:T633").Value = $arr
$ws.Range("D633").NumberFormat = "YYYY-MM-DD HH:MM:SS"
